$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above the current row 20 ("2021-10-15 / Primera"
# entry), shifting every existing row from 20 downward to 22 onward. This mirrors
# the two brand-new weekly price entries ("Extra" @44494 and "Primera" @44494)
# that the commit adds to the top of this block.
$ws.Rows.Item(20).Resize(2).Insert()

# New row 20: Espárragos, Extra, 2021-10-25 (serial 44494)
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = "2021-10-25"
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Extra"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 2000
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"

# New row 21: Espárragos, Primera, 2021-10-25 (serial 44494)
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = "2021-10-25"
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = "Espárragos"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1300
$ws.Range("M21").Value = 1300
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1300
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
